# Workshop_Presentation.pptx edit: "3 devient 3.5"
#
# On slide 7 ("3) Calcul des moyennes par catégorie"), the interpretation
# thresholds for the average score are updated:
#   "≤ 3 → À surveiller"  becomes  "≤ 3,5 → À surveiller"
#   "3 → Préoccupant"     becomes  "3,6 → Préoccupant" (now split across
#                                  two runs: "3,6 " + "→ Préoccupant")

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(7)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange

# --- Paragraph "≤ 3 → À surveiller" -> "≤ 3,5 → À surveiller" ---
$paraSurveiller = $tr.Paragraphs(8)
$paraSurveiller.Runs(1).Text = "≤ 3,5 → À surveiller"

# --- Paragraph "3 → Préoccupant" -> "3,6 " + "→ Préoccupant" ---
$paraPreoccupant = $tr.Paragraphs(9)

# Insert the new "3,6 " text in front of the existing run; this splits the
# paragraph into two runs: "3,6 " and the original "3 → Préoccupant".
$paraPreoccupant.Characters(1, 0).InsertBefore("3,6 ") | Out-Null

# Remove the leading "3 " (original run's first two characters, now shifted
# right by the 4 inserted characters) so the second run reads "→ Préoccupant".
$paraPreoccupant.Characters(5, 2).Text = ""
